# T460 - update foaie de parcurs (travel log) for B-151-VGT, iunie 2022, Alex Bora
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Km initiali (starting odometer reading)
$ws.Range("B12").Value = 399863

# Daily travel rows (Ziua / Km_parcursi / Locul deplasarii / Observatii utilizator)
$ws.Range("B15").Value = 85
$ws.Range("C15").Value = "Cluj-Apahida"
$ws.Range("D15").Value = "Interes Serviciu"

$ws.Range("B16").Value = 356
$ws.Range("C16").Value = "Cluj-Baia-Mare"
$ws.Range("D16").Value = "Interes Serviciu"

$ws.Range("B19").Value = 30
$ws.Range("C19").Value = "Acasa-Birou"
$ws.Range("D19").Value = " "

$ws.Range("B20").Value = 356
$ws.Range("C20").Value = "Cluj-Baia-Mare"
$ws.Range("D20").Value = "Interes Serviciu"

$ws.Range("B21").Value = 421
$ws.Range("C21").Value = "Cluj-Satu-Mare"
$ws.Range("D21").Value = "Interes Serviciu"

$ws.Range("B22").Value = 152
$ws.Range("C22").Value = "Cluj-Cmp. Turzii"
$ws.Range("D22").Value = "Interes Serviciu"

$ws.Range("B23").Value = 356
$ws.Range("C23").Value = "Cluj-Baia-Mare"
$ws.Range("D23").Value = "Interes Serviciu"

$ws.Range("B28").Value = 121
$ws.Range("C28").Value = "Cluj-Turda"
$ws.Range("D28").Value = "Interes Serviciu"

$ws.Range("B29").Value = 92
$ws.Range("C29").Value = "Cluj-Bontida"
$ws.Range("D29").Value = "Interes Serviciu"

$ws.Range("B30").Value = 156
$ws.Range("C30").Value = "Cluj-Zalau"
$ws.Range("D30").Value = "Interes Serviciu"

$ws.Range("B33").Value = 30
$ws.Range("C33").Value = "Acasa-Birou"
$ws.Range("D33").Value = " "

$ws.Range("B34").Value = 30
$ws.Range("C34").Value = "Acasa-Birou"
$ws.Range("D34").Value = " "

$ws.Range("B35").Value = 257
$ws.Range("C35").Value = "Cluj-Bistrita"
$ws.Range("D35").Value = "Interes Serviciu"

$ws.Range("B36").Value = 121
$ws.Range("C36").Value = "Cluj-Turda"
$ws.Range("D36").Value = "Interes Serviciu"

$ws.Range("B37").Value = 30
$ws.Range("C37").Value = "Acasa-Birou"
$ws.Range("D37").Value = " "

$ws.Range("B40").Value = 30
$ws.Range("C40").Value = "Acasa-Birou"
$ws.Range("D40").Value = " "

$ws.Range("B41").Value = 30
$ws.Range("C41").Value = "Acasa-Birou"
$ws.Range("D41").Value = " "

$ws.Range("B43").Value = 92
$ws.Range("C43").Value = "Cluj-Bontida"
$ws.Range("D43").Value = "Interes Serviciu"

# Totals
$ws.Range("B44").Value = 2876
$ws.Range("B45").Value = 402739
